$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-16.70343048998165, 1.791475864134279, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165),
    @(-16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165),
    @(-16.70343048998165, 2.170028642518187, 2.024905259129036, -16.70343048998165, 3.609865545262999, -16.70343048998165, 1.828716133148617, -16.70343048998165, 1.775058708092179, -16.70343048998165),
    @(-16.70343048998165, 1.582301333460042, -16.70343048998165, -16.70343048998165, -16.70343048998165, 2.842274971442194, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165),
    @(-16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165),
    @(2.534285972352167, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165),
    @(-16.70343048998165, -16.70343048998165, -16.70343048998165, 1.991094038700198, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165),
    @(3.828522977398892, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165),
    @(-16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, 2.277381288521799),
    @(-16.70343048998165, -16.70343048998165, -16.70343048998165, 2.834714857978484, -16.70343048998165, 3.124093643878027, -16.70343048998165, -16.70343048998165, -16.70343048998165, 1.986939269180629),
    @(-16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165),
    @(-16.70343048998165, -16.70343048998165, -16.70343048998165, 2.52463895045927, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, 1.440874244478844, 1.914772168691722),
    @(-16.70343048998165, -16.70343048998165, 1.128526147038312, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, 2.210198173612371),
    @(-16.70343048998165, -16.70343048998165, 0.9928646901835371, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165),
    @(-16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, 1.646926034709093, -16.70343048998165),
    @(-16.70343048998165, 2.382075683420311, 2.766857003074651, -16.70343048998165, -16.70343048998165, -16.70343048998165, 1.878465169056602, -16.70343048998165, 2.819934285708213, -16.70343048998165),
    @(-16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, -16.70343048998165, 1.775137705550981, 4.321915251887753, 1.875594433495853, -16.70343048998165),
    @(-16.70343048998165, -16.70343048998165, 1.403750820259324, -16.70343048998165, -16.70343048998165, -16.70343048998165, 1.806549357797753, -16.70343048998165, -16.70343048998165, -16.70343048998165),
    @(-16.70343048998165, 0.7834123859537577, 1.202643897364721, -16.70343048998165, 2.961789172414133, -16.70343048998165, 1.562939730095676, -16.70343048998165, -16.70343048998165, 1.480267883974925),
    @(-16.70343048998165, 1.076782450362457, -16.70343048998165, 1.648996635318762, -16.70343048998165, 2.039053298063764, 1.533821491493252, -16.70343048998165, -16.70343048998165, -16.70343048998165)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
